# Hjemme passive tweaks lichtwark deleted values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (subject-count header cells) - B1:E3 reworked with new passive data
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON)
$ws.Range("B2").Value = 10.72259672784131
$ws.Range("C2").Value = 15.430741855631982
$ws.Range("D2").Value = 8.5194839159582045
$ws.Range("E2").Value = 14.111444273894552

# Row 3 (STR)
$ws.Range("B3").Value = 12.959001009682119
$ws.Range("C3").Value = 15.711694743023926
$ws.Range("D3").Value = 15.496726544888144
$ws.Range("E3").Value = 15.341460935591842

# Selection now only spans the edited block
$ws.Range("B1:E3").Select()
